# "mejoras en el módulo de citas"
# Updates the cronograma (schedule) worksheet: bumps progress % on a few
# tasks, records a new "objetivo 100%" note, and appends a block of new
# planning notes below the existing summary table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Progress (% Avance) updates ---------------------------------------
$ws.Range("G7").Value = 0.9
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 0.4

# New note next to row 9 (merged H6:J6 area format is untouched here).
$ws.Range("H9").Value = "objetivo 100%"

# --- New planning notes --------------------------------------------------
# Order matters here: each first-use of a distinct string appends it to the
# shared-string table, so we write cells in the same order the strings were
# authored in.
$ws.Range("C26").Value = "gestionar fecha"
$ws.Range("C27").Value = "Gestionar Paciente"
$ws.Range("C28").Value = "Gestionar Doctor"
$ws.Range("C29").Value = "gestionar especialidad"
$ws.Range("C25").Value = "historia clínica completa"
$ws.Range("C23").Value = "voy a necesitar tiempo para investigar y/o preguntar como funciona todo y los documentos (los datos de la historia clínica, documento para tratamiento, historial tratamiento)"
$ws.Range("C30").Value = "Recetas"
$ws.Range("C31").Value = "voleta de venta"
$ws.Range("C33").Value = "el día miércoles te voy a dar un avance funcionando, le das el visto bueno y sigo agregando más cosas"
$ws.Range("C32").Value = "al ver que el sistema será más grande, propongo esto:"

# --- View state: scroll down and focus the last new note -----------------
$ws.Range("C31").Select()
$excel.ActiveWindow.ScrollRow = 10
